# Update "想去人数" (column F) figures on sheets "展览" and "全部类型"
# (and the matching rows shared with "本地生活"), reflecting newly
# generated output data.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value  = 21402
$wsExpo.Range("F3").Value  = 3288
$wsExpo.Range("F4").Value  = 855
$wsExpo.Range("F7").Value  = 800
$wsExpo.Range("F11").Value = 138
$wsExpo.Range("F12").Value = 563
$wsExpo.Range("F13").Value = 190
$wsExpo.Range("F14").Value = 339
$wsExpo.Range("F15").Value = 36
$wsExpo.Range("F17").Value = 172
$wsExpo.Range("F20").Value = 78
$wsExpo.Range("F21").Value = 148

# --- Sheet: 本地生活 ---
$wsLocal = $wb.Worksheets.Item("本地生活")
$wsLocal.Range("F2").Value = 6165
$wsLocal.Range("F4").Value = 727
$wsLocal.Range("F5").Value = 1727
$wsLocal.Range("F6").Value = 82

# --- Sheet: 全部类型 ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value  = 6165
$wsAll.Range("F4").Value  = 727
$wsAll.Range("F5").Value  = 1727
$wsAll.Range("F6").Value  = 21402
$wsAll.Range("F7").Value  = 3288
$wsAll.Range("F8").Value  = 855
$wsAll.Range("F10").Value = 82
$wsAll.Range("F13").Value = 800
$wsAll.Range("F20").Value = 138
$wsAll.Range("F23").Value = 563
$wsAll.Range("F25").Value = 190
$wsAll.Range("F27").Value = 339
$wsAll.Range("F29").Value = 36
$wsAll.Range("F32").Value = 172
$wsAll.Range("F37").Value = 78
$wsAll.Range("F43").Value = 148
